# Update the questionnaire worksheet with the new, larger set of sections
# and questions (Waste details / Input materials or fuels / Emissions to air).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8 previously was the last row in the table and had no "options" (D)
# cell at all. Give it the same blank-but-styled D cell the other existing
# data rows (4-7) have, by copying formatting across from a neighboring
# cell on the same row.
$ws.Range("E8").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D8").Value = ""

# Rows 9-33 are brand new. Stamp them with the same look as the existing
# data rows (style "2" on columns A/B/C/E, no D cell at all) by copying
# formatting down from row 8 before writing any values.
$ws.Range("A8:C8").Copy()
$ws.Range("A9:C33").PasteSpecial(-4122)
$ws.Range("E8").Copy()
$ws.Range("E9:E33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 4: "Waste details" section / "waste composition" question.
$ws.Range("A4").Value = "Waste details"
$ws.Range("B4").Value = "waste composition"
$ws.Range("C4").Value = "text"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "yes"

# "Input materials or fuels" section questions (rows 5-14).
$inputMaterials = @(
    "Waste treated quantity",
    "Sodium bicarbonate used",
    "Activated carbon",
    "Ammonia",
    "Sodium hypochlorite (NaClO)",
    "Scale and corrosion inhibitors",
    "Alkalinizing agent",
    "Deoxidizer",
    "Water",
    "Auxiliary fuel"
)

$row = 5
foreach ($q in $inputMaterials) {
    $ws.Cells.Item($row, 1).Value = "Input materials or fuels"
    $ws.Cells.Item($row, 2).Value = $q
    $ws.Cells.Item($row, 3).Value = "number"
    $ws.Cells.Item($row, 5).Value = "yes"
    $row++
}

# "Emissions to air" section questions (rows 15-33).
$emissions = @(
    "Total Carbone dioxide",
    "Hydrochloric acid (HCl)",
    "Hydrofluoric acid (HF)",
    "Sulfur dioxide (SO2)",
    "Nitrogen dioxide (NO2)",
    "Ammonia (NH3)",
    "Mercury (Hg)",
    "Zinc (Zn)",
    "Dioxins and furans (PCDD/PCDF)",
    "Carbon monoxide (CO)",
    "Total dust (< 10 µm)",
    "Polycyclic aromatic hydrocarbons (PAHs)",
    "Cadmium (Cd)",
    "Thallium (Tl)",
    "Antimony (Sb)",
    "Arsenic (As)",
    "Lead (Pb)",
    "Chromium (Cr)",
    "Cobalt (Co)"
)

foreach ($q in $emissions) {
    $ws.Cells.Item($row, 1).Value = "Emissions to air"
    $ws.Cells.Item($row, 2).Value = $q
    $ws.Cells.Item($row, 3).Value = "number"
    $ws.Cells.Item($row, 5).Value = "yes"
    $row++
}

# Move the active selection to reflect the final cursor position after editing.
$ws.Range("H15").Select()
